# Apply the "add monte_carlo and update database" update to the
# quarterly rial_cumulative income statement workbook.
#
# The sheet "Overview" holds a cumulative (year-to-date) income statement
# whose last reporting period is the column M ("12 ماهه منتهی به 1401/12").
# This update re-publishes that period's figures under a later publish
# date and refreshes the cumulative numbers that changed as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Publish-date labels (row 9) ---------------------------------------
# Column I (3 ماهه منتهی به 1401/09 column pairing) publish date label.
$ws.Range("I9").Value = "1402-03-13 (10)"
# Column M (12 ماهه منتهی به 1401/12) publish date label.
$ws.Range("M9").Value = "1402-03-13 (2)"

# --- Updated cumulative figures (column M) ------------------------------
$ws.Range("M14").Value = -9451842
$ws.Range("M17").Value = 12883983
$ws.Range("M18").Value = -7820256
$ws.Range("M20").Value = 32464667
$ws.Range("M21").Value = -742432
$ws.Range("M22").Value = 31722235
$ws.Range("M24").Value = 31722235
$ws.Range("M25").Value = 1025
$ws.Range("M27").Value = 1025
